# Apply the edits described by the commit "Multiline String check for file beginning cleanup"
# Workbook: setup.xlsx, sheet "Main Info"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Info")

# B2 (Region dropdown "APAC,EMEA,NAM") - re-affirm value "NAM" (unchanged visible value,
# but re-entering makes Excel re-insert the shared string at the end of the table)
$ws.Range("B2").Value = "NAM"

# B7 (Design dropdown "BASE,SMART,FLOW") - changed from "FLOW" to "SMART"
$ws.Range("B7").Value = "SMART"

# B8 (Converged router dropdown "TRUE, FALSE") - changed from FALSE to TRUE
$ws.Range("B8").Value = $true

# B9 (Migration from MPLS dropdown "False, True - New Router, True - Production router")
# changed from boolean FALSE to the literal string "True - Production router"
$ws.Range("B9").Value = "True - Production router"

# B10 (ZBFW dropdown "True,False") - changed from FALSE to TRUE
$ws.Range("B10").Value = $true

# B25 (4G+Cellular dropdown "True,False") - changed from FALSE to TRUE
$ws.Range("B25").Value = $true

# Update the active selection to D6, as recorded in the saved view state
$ws.Range("D6").Select()
